# daily auto push: 2026-02-01 13:51 UTC
#
# Inserts a new data row at row 763 (shifting all following rows down by
# one) containing a new observation for 2026/02/01 (日, 19時, ランキング166).
# This mirrors the upstream diff, which pushed one extra row into the
# middle of the existing "2026/12/29 .. 2027/01/05" block and bumped the
# sheet dimension from A1:D804 to A1:D805.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$targetRow = 763

# Shift rows 763..804 down to 764..805, creating a blank row 763.
$ws.Rows.Item($targetRow).Insert()

# Column A holds a date formatted as plain text (e.g. "2026/12/29"), not a
# real Excel date. Excel's COM layer auto-converts date-looking text typed
# into a General formatted cell into a serial date, so the cell is briefly
# switched to Text format ("@") while the literal string is written, then
# restored to the default "Normal" style so no stray formatting is left
# behind (matching the unstyled cells used throughout this sheet).
$dateCell = $ws.Cells.Item($targetRow, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2026/02/01"
$dateCell.Style = "Normal"

$ws.Cells.Item($targetRow, 2).Value = "日"
$ws.Cells.Item($targetRow, 3).Value = 19
$ws.Cells.Item($targetRow, 4).Value = 166
